$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 4 (pushing all existing data, e.g. old row 4, down to row 6, etc.)
$ws.Rows("4:5").Insert()

# New row 4: Jengibre, Primera, new weekly price entry (fecha 44515)
$ws.Range("A4").Value2 = 9
$ws.Range("B4").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C4").Value2 = "Metropolitana"
$ws.Range("D4").Value2 = 44515
$ws.Range("E4").Value2 = 13
$ws.Range("F4").Value2 = 100114007
$ws.Range("G4").Value2 = "Jengibre"
$ws.Range("H4").Value2 = "Sin especificar"
$ws.Range("I4").Value2 = "Primera"
$ws.Range("J4").Value2 = 1060
$ws.Range("K4").Value2 = 16000
$ws.Range("L4").Value2 = 18000
$ws.Range("M4").Value2 = 17000
$ws.Range("N4").Value2 = "$/caja 13 kilos"
$ws.Range("O4").Value2 = "Perú"
$ws.Range("P4").Value2 = 1308
$ws.Range("Q4").Value2 = 13
$ws.Range("R4").Value2 = "Hortaliza"

# New row 5: Jengibre, Segunda, new weekly price entry (fecha 44515)
$ws.Range("A5").Value2 = 9
$ws.Range("B5").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C5").Value2 = "Metropolitana"
$ws.Range("D5").Value2 = 44515
$ws.Range("E5").Value2 = 13
$ws.Range("F5").Value2 = 100114007
$ws.Range("G5").Value2 = "Jengibre"
$ws.Range("H5").Value2 = "Sin especificar"
$ws.Range("I5").Value2 = "Segunda"
$ws.Range("J5").Value2 = 610
$ws.Range("K5").Value2 = 14000
$ws.Range("L5").Value2 = 14000
$ws.Range("M5").Value2 = 14000
$ws.Range("N5").Value2 = "$/caja 13 kilos"
$ws.Range("O5").Value2 = "Perú"
$ws.Range("P5").Value2 = 1077
$ws.Range("Q5").Value2 = 13
$ws.Range("R5").Value2 = "Hortaliza"

# Apply the same date number format used by other date cells in column D
$ws.Range("D4").NumberFormat = $ws.Range("D6").NumberFormat
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat
